$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 6; existing rows 6-283 shift down to 7-284.
$ws.Rows.Item(6).Insert()

# The new blank row 6 should carry the same "entity" data as the row now
# sitting at row 7 (all rows share identical Mercado/Categoria/etc.), so
# copy row 7 into row 6, then set the new week's date on D6.
$ws.Range("A7:R7").Copy()
$ws.Range("A6:R6").PasteSpecial()

$ws.Range("D6").Value = 44631
